$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (originally row 26) and the "SC 92" row
# (originally row 28, which becomes row 27 after the first deletion).
# Everything below shifts up, turning the 34-data-row table into a
# 32-data-row table (A1:F33).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the individual cell value changes (re-sampled "missing data" cells).
$ws.Range("D2").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("D5").Value = -14.4
$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("C8").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = -6.1
$ws.Range("E11").Value = -7.9
$ws.Range("C12").Value = 12.5
$ws.Range("E12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("E14").Value = -5.4
$ws.Range("E16").Value = ""
$ws.Range("C17").Value = 11.2
$ws.Range("E17").Value = ""
$ws.Range("C18").Value = 11.5
$ws.Range("C19").Value = ""
$ws.Range("E19").Value = -6.5
$ws.Range("C20").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E22").Value = -6.1
$ws.Range("C23").Value = 12.2
$ws.Range("D24").Value = -13.9
$ws.Range("E25").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("D30").Value = -13.6
$ws.Range("E31").Value = -8.1
$ws.Range("B32").Value = ""
